$d = $word.ActiveDocument

# 1. Replace the repeated "verser" with "servir" in the target sentence,
#    to avoid repeating "verser" twice in a row.
$d.Content.Find.Execute("préféra se verser de l’eau", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "préféra se servir de l’eau", 2)

# 2. Find the trailing lone space run after "... à la place." and replace it
#    with the new sentences that continue the scene.
$d.Content.Find.Execute("à la place. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "à la place. L’eau était fraiche. Cela lui fit du bien. Elle s’habilla d’une robe de lin légère, lui laissant une grande liberté de mouvement. Elle ne voulait pas de couper la respiration avec un corsage trop serré, et enfila des sandales plates en cuir. Une fois prête, elle sortit sans bruit de ses appartements, vérifiant que son époux était endormi quand elle ferma la porte. A quelques pas de la porte, deux hommes montaient la garde.", 2)
